$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.677.24"
$ws.Range("E2").Value = "  +3.38%  "

$ws.Range("D3").Value = "3.356.67"
$ws.Range("E3").Value = "  +4.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "193.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "592.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "

$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("D12").Value = "3.939.12"
$ws.Range("E12").Value = "  +4.25%  "

$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.15%  "

$ws.Range("D15").Value = "69.668.72"
$ws.Range("E15").Value = "  +3.27%  "

$ws.Range("E16").Value = "  +2.19%  "

$ws.Range("D17").Value = "3.343.76"
$ws.Range("E17").Value = "  +4.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "444.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.61%  "

$ws.Range("E22").Value = "  +4.00%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "3.492.45"
$ws.Range("E24").Value = "  +4.04%  "

$ws.Range("E25").Value = "  +4.14%  "

$ws.Range("E26").Value = "  +1.08%  "

$ws.Range("E27").Value = "  +4.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("E33").Value = "  +3.32%  "

$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "

$ws.Range("E37").Value = "  +3.16%  "

$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("E40").Value = "  +1.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "

$ws.Range("D42").Value = "2.755.87"
$ws.Range("E42").Value = "  +5.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("E44").Value = "  +3.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "345.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.58%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0688"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("E49").Value = "  +3.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.24%  "

$ws.Range("E51").Value = "  +4.62%  "

